# Weekly update: insert a new row of data at the top of the Haba price
# series (row 8, just below the fixed header block in rows 1-7) and push
# every existing data row down by one. The oldest row (previously row 104)
# ends up as the new row 105, with the sheet's used range growing from
# A1:R104 to A1:R105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..104 down to 9..105, duplicating row 8's formatting
# (including the date-format style on column D) into the newly freed row 8.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with this week's record.
$ws.Cells.Item(8, 1).Value = 3                                  # A8  Mercado ID
$ws.Cells.Item(8, 2).Value = "Femacal de La Calera"              # B8  Mercado
$ws.Cells.Item(8, 3).Value = "Coquimbo"                          # C8  Región
$ws.Cells.Item(8, 4).Value = 44530                                # D8  Fecha
$ws.Cells.Item(8, 5).Value = 5                                   # E8  Codreg
$ws.Cells.Item(8, 6).Value = 100112026                            # F8  Categoría ID
$ws.Cells.Item(8, 7).Value = "Haba"                               # G8  Categoría
$ws.Cells.Item(8, 8).Value = "Sin especificar"                   # H8  Variedad
$ws.Cells.Item(8, 9).Value = "Primera"                           # I8  Calidad
$ws.Cells.Item(8, 10).Value = 85                                  # J8  Volumen
$ws.Cells.Item(8, 11).Value = 7500                                # K8  Precio mínimo
$ws.Cells.Item(8, 12).Value = 8000                                # L8  Precio máximo
$ws.Cells.Item(8, 13).Value = 7735                                # M8  Precio promedio ponderado
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"                  # N8  Unidad de comercialización
$ws.Cells.Item(8, 15).Value = "Provincia de Quillota"            # O8  Origen
$ws.Cells.Item(8, 16).Value = 309                                 # P8  Precio $/Kg
$ws.Cells.Item(8, 17).Value = 25                                  # Q8  Kg o Unidades
$ws.Cells.Item(8, 18).Value = "Hortaliza"                         # R8  Clasificación
